$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.141.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "'3.542.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'598.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.65%  "
$ws.Range("D6").Value = "'138.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("D7").Value = "'3.541.85"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.22%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'0.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("D11").Value = "'6.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("E12").Value = "  +4.23%  "
$ws.Range("D13").Value = "'4.145.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.33%  "
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").Value = "'27.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.15%  "
$ws.Range("D16").Value = "'3.548.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "'65.083.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'10.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.20%  "
$ws.Range("D20").Value = "'5.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "'14.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.52%  "
$ws.Range("D22").Value = "'392.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("E23").Value = "  +4.67%  "
$ws.Range("D24").Value = "'3.684.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.32%  "
$ws.Range("D25").Value = "'73.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D28").Value = "'7.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.71%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("D32").Value = "'3.562.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.57%  "
$ws.Range("D33").Value = "'1.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +22.04%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'23.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "'1.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.57%  "
$ws.Range("D38").Value = "'6.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.42%  "
$ws.Range("D39").Value = "'168.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  +8.48%  "
$ws.Range("D41").Value = "'0.0807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.48%  "
$ws.Range("D42").Value = "'0.825"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("D43").Value = "'26.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +22.48%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("E47").Value = "  +10.38%  "
$ws.Range("E48").Value = "  +5.94%  "
$ws.Range("E49").Value = "  +6.34%  "
$ws.Range("D50").Value = "'2.403.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.62%  "
$ws.Range("D51").Value = "'312.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.39%  "
